$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as text in the source sheet.
# Force text number format per-cell before assigning so that numeric-looking
# values (e.g. "3.00", "0.990") are preserved exactly as text, not coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.913.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.551.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.771.37"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.551.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.909.58"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.22"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.12"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.95"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.411.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.521"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.990"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.49"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.686.12"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.21"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0959"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.04%  "
